# This script reorders the weekly price observation rows (rows 2-30, columns A-T)
# of the active worksheet according to a fixed row permutation, as produced by the
# target diff. Column A1:T1 (headers) are left untouched.
#
# Mapping is expressed as: for each destination row (2..30), which source row
# (2..30) supplies its data in the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (both are absolute worksheet row numbers)
$rowMap = @{
    2  = 12
    3  = 22
    4  = 4
    5  = 15
    6  = 16
    7  = 13
    8  = 14
    9  = 24
    10 = 25
    11 = 2
    12 = 3
    13 = 7
    14 = 8
    15 = 26
    16 = 19
    17 = 20
    18 = 21
    19 = 30
    20 = 5
    21 = 17
    22 = 18
    23 = 27
    24 = 28
    25 = 29
    26 = 23
    27 = 10
    28 = 11
    29 = 9
    30 = 6
}

$firstRow = 2
$lastRow = 30
$firstCol = 1   # A
$lastCol = 20   # T

# Read the entire current block (values only, fast/reliable for dates & numbers)
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$srcData = $srcRange.Value2

$rowCount = $lastRow - $firstRow + 1
$colCount = $lastCol - $firstCol + 1

# Build the new block in the desired (permuted) row order.
# Note: arrays created via New-Object are 0-based, while the Value2 array
# returned from COM (srcData) is 1-based; indices below account for that.
$newData = New-Object 'object[,]' $rowCount, $colCount
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $destIdx = $destRow - $firstRow        # 0-based
    $srcIdx = $srcRow - $firstRow + 1      # 1-based
    for ($c = 1; $c -le $colCount; $c++) {
        $newData[$destIdx, ($c - 1)] = $srcData[$srcIdx, $c]
    }
}

# Write the permuted block back in a single pass
$destRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$destRange.Value2 = $newData
